# Weekly update: insert a new daily record as row 118 ("Hortaliza, Feria
# Lagunitas de Puerto Montt - Cilantro"), pushing the existing rows 118:224
# down to 119:225 (dimension grows from A1:R224 to A1:R225).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 118; this shifts rows 118-224
# down to 119-225 and copies formatting (e.g. the date style on column D)
# from the row above, matching native Excel "Insert Row" behavior.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new day's reading.
$ws.Range("A118").Value = 4
$ws.Range("B118").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C118").Value = 'Los Lagos'
$ws.Range("D118").Value = 44587
$ws.Range("E118").Value = 10
$ws.Range("F118").Value = 100112040
$ws.Range("G118").Value = 'Cilantro'
$ws.Range("H118").Value = 'Sin especificar'
$ws.Range("I118").Value = 'Primera'
$ws.Range("J118").Value = 20
$ws.Range("K118").Value = 10000
$ws.Range("L118").Value = 10000
$ws.Range("M118").Value = 10000
$ws.Range("N118").Value = '$/docena de atados (2 kilos)'
$ws.Range("O118").Value = 'Región de La Araucanía'
$ws.Range("P118").Value = 5000
$ws.Range("Q118").Value = 2
$ws.Range("R118").Value = 'Hortaliza'
